$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was incremented
# by one day (46061 -> 46062) for every data row (rows 2 through 53).
for ($r = 2; $r -le 53; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2
    $cell.Value = $current + 1
}
